$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 259.8889
$ws.Cells.Item(2, 9).Value = 219.85715
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 219.85715
$ws.Cells.Item(2, 12).Value = 400
$ws.Cells.Item(2, 13).Value = -106.85715
$ws.Cells.Item(2, 14).Value = -626
$ws.Cells.Item(28, 8).Value = 1098.8889
$ws.Cells.Item(28, 9).Value = 1550
$ws.Cells.Item(28, 10).Value = 873.3333
$ws.Cells.Item(28, 11).Value = 1550
$ws.Cells.Item(28, 12).Value = 873.3333
$ws.Cells.Item(28, 13).Value = -1065
$ws.Cells.Item(28, 14).Value = -1843.3333
$ws.Cells.Item(38, 8).Value = 2928
$ws.Cells.Item(38, 9).Value = 276
$ws.Cells.Item(38, 10).Value = 10000
$ws.Cells.Item(38, 11).Value = 828
$ws.Cells.Item(38, 12).Value = 30000
$ws.Cells.Item(38, 13).Value = -456
$ws.Cells.Item(38, 14).Value = -30744
$ws.Cells.Item(40, 8).Value = 4584.1177
$ws.Cells.Item(40, 9).Value = 3413.4
$ws.Cells.Item(40, 11).Value = 3413.4
$ws.Cells.Item(40, 13).Value = -3238.4
$ws.Cells.Item(41, 8).Value = 227.4
$ws.Cells.Item(41, 9).Value = 209.5
$ws.Cells.Item(41, 11).Value = 209.5
$ws.Cells.Item(41, 13).Value = 230.5
$ws.Cells.Item(74, 8).Value = 6505.533
$ws.Cells.Item(74, 9).Value = 4816.6
$ws.Cells.Item(74, 11).Value = 4816.6
$ws.Cells.Item(74, 13).Value = -3880.6
$ws.Cells.Item(77, 8).Value = 6505.533
$ws.Cells.Item(77, 9).Value = 4816.6
$ws.Cells.Item(77, 11).Value = 24083
$ws.Cells.Item(77, 13).Value = -19403
$ws.Cells.Item(88, 8).Value = 2449.4443
$ws.Cells.Item(88, 9).Value = 1962.2858
$ws.Cells.Item(88, 10).Value = 2759.4546
$ws.Cells.Item(88, 11).Value = 1962.2858
$ws.Cells.Item(88, 12).Value = 2759.4546
$ws.Cells.Item(88, 13).Value = -1556.2858
$ws.Cells.Item(88, 14).Value = -3571.4546
$ws.Cells.Item(91, 8).Value = 2449.4443
$ws.Cells.Item(91, 9).Value = 1962.2858
$ws.Cells.Item(91, 10).Value = 2759.4546
$ws.Cells.Item(91, 11).Value = 1962.2858
$ws.Cells.Item(91, 12).Value = 2759.4546
$ws.Cells.Item(91, 13).Value = -558.2858000000001
$ws.Cells.Item(91, 14).Value = -5567.4546
$ws.Cells.Item(96, 8).Value = 405.72726
$ws.Cells.Item(96, 9).Value = 436.4
$ws.Cells.Item(96, 10).Value = 99
$ws.Cells.Item(96, 11).Value = 1309.2
$ws.Cells.Item(96, 12).Value = 297
$ws.Cells.Item(96, 13).Value = 63.80000000000018
$ws.Cells.Item(96, 14).Value = -3043
$ws.Cells.Item(99, 8).Value = 188.875
$ws.Cells.Item(99, 9).Value = 180.14285
$ws.Cells.Item(99, 10).Value = 250
$ws.Cells.Item(99, 11).Value = 540.4285500000001
$ws.Cells.Item(99, 12).Value = 750
$ws.Cells.Item(99, 13).Value = 957.5714499999999
$ws.Cells.Item(99, 14).Value = -3746
$ws.Cells.Item(107, 8).Value = 541.8570999999999
$ws.Cells.Item(107, 9).Value = 585.76
$ws.Cells.Item(107, 11).Value = 585.76
$ws.Cells.Item(107, 13).Value = 1334.24
$ws.Cells.Item(111, 8).Value = 779.8
$ws.Cells.Item(111, 9).Value = 349.5
$ws.Cells.Item(111, 10).Value = 1066.6666
$ws.Cells.Item(111, 11).Value = 1048.5
$ws.Cells.Item(111, 12).Value = 3199.9998
$ws.Cells.Item(111, 13).Value = 2018.5
$ws.Cells.Item(111, 14).Value = -9333.9998
$ws.Cells.Item(132, 8).Value = 4679.5312
$ws.Cells.Item(132, 9).Value = 1096.1666
$ws.Cells.Item(132, 11).Value = 3288.4998
$ws.Cells.Item(132, 13).Value = -758.4998000000001
$ws.Cells.Item(133, 8).Value = 74749.5
$ws.Cells.Item(133, 10).Value = 74749.5
$ws.Cells.Item(133, 12).Value = 74749.5
$ws.Cells.Item(133, 14).Value = -84869.5
$ws.Cells.Item(137, 8).Value = 1473.5555
$ws.Cells.Item(137, 9).Value = 1126.6666
$ws.Cells.Item(137, 11).Value = 3379.9998
$ws.Cells.Item(137, 13).Value = -829.9998000000001
$ws.Cells.Item(141, 8).Value = 3207.2917
$ws.Cells.Item(141, 9).Value = 2973.95
$ws.Cells.Item(141, 11).Value = 8921.849999999999
$ws.Cells.Item(141, 13).Value = -3741.849999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 1505
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 14).Value = $null
$ws.Cells.Item(21, 8).Value = 7330
$ws.Cells.Item(21, 9).Value = 7330
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 7330
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = -6956
$ws.Cells.Item(21, 14).Value = $null
$ws.Cells.Item(32, 8).Value = 1256.0137
$ws.Cells.Item(32, 9).Value = 1256.0137
$ws.Cells.Item(32, 11).Value = 1256.0137
$ws.Cells.Item(32, 13).Value = -969.0137
$ws.Cells.Item(37, 8).Value = 11250
$ws.Cells.Item(37, 9).Value = 1500
$ws.Cells.Item(37, 11).Value = 1500
$ws.Cells.Item(37, 13).Value = -1227
$ws.Cells.Item(61, 8).Value = 10172.381
$ws.Cells.Item(61, 9).Value = 11163.462
$ws.Cells.Item(61, 11).Value = 11163.462
$ws.Cells.Item(61, 13).Value = -10951.462
$ws.Cells.Item(63, 8).Value = 5073
$ws.Cells.Item(63, 9).Value = 3666.6667
$ws.Cells.Item(63, 11).Value = 3666.6667
$ws.Cells.Item(63, 13).Value = -2980.6667
$ws.Cells.Item(66, 8).Value = 5073
$ws.Cells.Item(66, 9).Value = 3666.6667
$ws.Cells.Item(66, 11).Value = 18333.3335
$ws.Cells.Item(66, 13).Value = -14901.3335
$ws.Cells.Item(74, 8).Value = 3576.4211
$ws.Cells.Item(74, 9).Value = 2977.5
$ws.Cells.Item(74, 11).Value = 2977.5
$ws.Cells.Item(74, 13).Value = -2103.5
$ws.Cells.Item(77, 8).Value = 3576.4211
$ws.Cells.Item(77, 9).Value = 2977.5
$ws.Cells.Item(77, 11).Value = 14887.5
$ws.Cells.Item(77, 13).Value = -10519.5
$ws.Cells.Item(110, 8).Value = 5038.5884
$ws.Cells.Item(110, 9).Value = 3201.5
$ws.Cells.Item(110, 11).Value = 3201.5
$ws.Cells.Item(110, 13).Value = -1156.5
$ws.Cells.Item(122, 8).Value = 2754.2
$ws.Cells.Item(122, 9).Value = 1906.1428
$ws.Cells.Item(122, 10).Value = 4733
$ws.Cells.Item(122, 11).Value = 5718.428400000001
$ws.Cells.Item(122, 12).Value = 14199
$ws.Cells.Item(122, 13).Value = -3268.428400000001
$ws.Cells.Item(122, 14).Value = -19099
$ws.Cells.Item(136, 8).Value = 10172.381
$ws.Cells.Item(136, 9).Value = 11163.462
$ws.Cells.Item(136, 11).Value = 33490.386
$ws.Cells.Item(136, 13).Value = -30940.386

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 11999
$ws.Cells.Item(26, 9).Value = 11999
$ws.Cells.Item(26, 11).Value = 11999
$ws.Cells.Item(26, 13).Value = -11707
$ws.Cells.Item(57, 8).Value = 48992.75
$ws.Cells.Item(57, 10).Value = 48992.75
$ws.Cells.Item(57, 12).Value = 48992.75
$ws.Cells.Item(57, 14).Value = -50432.75
$ws.Cells.Item(86, 8).Value = 4333
$ws.Cells.Item(86, 9).Value = 3999.5
$ws.Cells.Item(86, 11).Value = 3999.5
$ws.Cells.Item(86, 13).Value = -2876.5
$ws.Cells.Item(89, 8).Value = 4333
$ws.Cells.Item(89, 9).Value = 3999.5
$ws.Cells.Item(89, 11).Value = 19997.5
$ws.Cells.Item(89, 13).Value = -14381.5
$ws.Cells.Item(94, 8).Value = 2232.9092
$ws.Cells.Item(94, 9).Value = 1944.6428
$ws.Cells.Item(94, 10).Value = 2737.375
$ws.Cells.Item(94, 11).Value = 1944.6428
$ws.Cells.Item(94, 12).Value = 2737.375
$ws.Cells.Item(94, 13).Value = -1493.6428
$ws.Cells.Item(94, 14).Value = -3639.375
$ws.Cells.Item(105, 8).Value = 937
$ws.Cells.Item(105, 9).Value = 937
$ws.Cells.Item(105, 11).Value = 937
$ws.Cells.Item(105, 13).Value = 810
$ws.Cells.Item(107, 8).Value = 1757.1765
$ws.Cells.Item(107, 9).Value = 1248.0714
$ws.Cells.Item(107, 11).Value = 1248.0714
$ws.Cells.Item(107, 13).Value = 671.9286
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 14).Value = $null
$ws.Cells.Item(134, 8).Value = 3704.5908
$ws.Cells.Item(134, 9).Value = 3704.5908
$ws.Cells.Item(134, 11).Value = 11113.7724
$ws.Cells.Item(134, 13).Value = -8578.7724
$ws.Cells.Item(136, 8).Value = 48992.75
$ws.Cells.Item(136, 10).Value = 48992.75
$ws.Cells.Item(136, 12).Value = 48992.75
$ws.Cells.Item(136, 14).Value = -59192.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = $null
$ws.Cells.Item(31, 8).Value = 4437.4707
$ws.Cells.Item(31, 9).Value = 2311.1538
$ws.Cells.Item(31, 11).Value = 2311.1538
$ws.Cells.Item(31, 13).Value = -2016.1538
$ws.Cells.Item(34, 8).Value = 4437.4707
$ws.Cells.Item(34, 9).Value = 2311.1538
$ws.Cells.Item(34, 11).Value = 2311.1538
$ws.Cells.Item(34, 13).Value = -2109.1538
$ws.Cells.Item(58, 8).Value = 5709.36
$ws.Cells.Item(58, 9).Value = 4189.154
$ws.Cells.Item(58, 10).Value = 7356.25
$ws.Cells.Item(58, 11).Value = 4189.154
$ws.Cells.Item(58, 12).Value = 7356.25
$ws.Cells.Item(58, 13).Value = -3986.154
$ws.Cells.Item(58, 14).Value = -7762.25
$ws.Cells.Item(62, 8).Value = 6643.143
$ws.Cells.Item(62, 9).Value = 6750.3335
$ws.Cells.Item(62, 11).Value = 6750.3335
$ws.Cells.Item(62, 13).Value = -6126.3335
$ws.Cells.Item(65, 8).Value = 6643.143
$ws.Cells.Item(65, 9).Value = 6750.3335
$ws.Cells.Item(65, 11).Value = 33751.6675
$ws.Cells.Item(65, 13).Value = -30631.6675
$ws.Cells.Item(80, 8).Value = 37141.453
$ws.Cells.Item(80, 10).Value = 37141.453
$ws.Cells.Item(80, 12).Value = 37141.453
$ws.Cells.Item(80, 14).Value = -39387.453
$ws.Cells.Item(83, 8).Value = 37141.453
$ws.Cells.Item(83, 10).Value = 37141.453
$ws.Cells.Item(83, 12).Value = 111424.359
$ws.Cells.Item(83, 14).Value = -122656.359
$ws.Cells.Item(99, 8).Value = 4329.8
$ws.Cells.Item(99, 9).Value = 4329.8
$ws.Cells.Item(99, 11).Value = 4329.8
$ws.Cells.Item(99, 13).Value = -2831.8
$ws.Cells.Item(120, 8).Value = 39441.832
$ws.Cells.Item(120, 9).Value = 39325
$ws.Cells.Item(120, 10).Value = 39465.2
$ws.Cells.Item(120, 11).Value = 39325
$ws.Cells.Item(120, 12).Value = 39465.2
$ws.Cells.Item(120, 13).Value = -35696
$ws.Cells.Item(120, 14).Value = -46723.2
$ws.Cells.Item(121, 8).Value = 187965.2
$ws.Cells.Item(121, 10).Value = 187965.2
$ws.Cells.Item(121, 12).Value = 187965.2
$ws.Cells.Item(121, 14).Value = -190585.2
$ws.Cells.Item(122, 8).Value = 3067.4783
$ws.Cells.Item(122, 9).Value = 2525.9
$ws.Cells.Item(122, 10).Value = 6678
$ws.Cells.Item(122, 11).Value = 7577.700000000001
$ws.Cells.Item(122, 12).Value = 20034
$ws.Cells.Item(122, 13).Value = -5127.700000000001
$ws.Cells.Item(122, 14).Value = -24934
$ws.Cells.Item(126, 8).Value = 4329.8
$ws.Cells.Item(126, 9).Value = 4329.8
$ws.Cells.Item(126, 11).Value = 12989.4
$ws.Cells.Item(126, 13).Value = -10519.4
$ws.Cells.Item(133, 8).Value = 52013
$ws.Cells.Item(133, 10).Value = 55769.5
$ws.Cells.Item(133, 12).Value = 55769.5
$ws.Cells.Item(133, 14).Value = -60829.5
$ws.Cells.Item(136, 8).Value = 5709.36
$ws.Cells.Item(136, 9).Value = 4189.154
$ws.Cells.Item(136, 10).Value = 7356.25
$ws.Cells.Item(136, 11).Value = 12567.462
$ws.Cells.Item(136, 12).Value = 22068.75
$ws.Cells.Item(136, 13).Value = -10017.462
$ws.Cells.Item(136, 14).Value = -27168.75
$ws.Cells.Item(138, 8).Value = 49455.6
$ws.Cells.Item(138, 10).Value = 159294.33
$ws.Cells.Item(138, 12).Value = 159294.33
$ws.Cells.Item(138, 14).Value = -169574.33
$ws.Cells.Item(140, 8).Value = 127496.25
$ws.Cells.Item(140, 10).Value = 127496.25
$ws.Cells.Item(140, 12).Value = 127496.25
$ws.Cells.Item(140, 14).Value = -137856.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 13622
$ws.Cells.Item(3, 9).Value = 9795.200000000001
$ws.Cells.Item(3, 11).Value = 29385.6
$ws.Cells.Item(3, 13).Value = -29273.6
$ws.Cells.Item(113, 8).Value = 665.9
$ws.Cells.Item(113, 10).Value = 680.5
$ws.Cells.Item(113, 12).Value = 2041.5
$ws.Cells.Item(113, 14).Value = -6381.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 12794.75
$ws.Cells.Item(2, 9).Value = 226.5
$ws.Cells.Item(2, 10).Value = 50499.5
$ws.Cells.Item(2, 11).Value = 226.5
$ws.Cells.Item(2, 12).Value = 50499.5
$ws.Cells.Item(2, 13).Value = -113.5
$ws.Cells.Item(2, 14).Value = -50725.5
$ws.Cells.Item(80, 8).Value = 4525.6113
$ws.Cells.Item(80, 9).Value = 4035.3635
$ws.Cells.Item(80, 10).Value = 5296
$ws.Cells.Item(80, 11).Value = 4035.3635
$ws.Cells.Item(80, 12).Value = 5296
$ws.Cells.Item(80, 13).Value = -3037.3635
$ws.Cells.Item(80, 14).Value = -7292
$ws.Cells.Item(83, 8).Value = 4525.6113
$ws.Cells.Item(83, 9).Value = 4035.3635
$ws.Cells.Item(83, 10).Value = 5296
$ws.Cells.Item(83, 11).Value = 20176.8175
$ws.Cells.Item(83, 12).Value = 26480
$ws.Cells.Item(83, 13).Value = -15184.8175
$ws.Cells.Item(83, 14).Value = -36464
$ws.Cells.Item(102, 8).Value = 2711.389
$ws.Cells.Item(102, 9).Value = 2144.0625
$ws.Cells.Item(102, 11).Value = 2144.0625
$ws.Cells.Item(102, 13).Value = -522.0625
$ws.Cells.Item(113, 8).Value = 445956.66
$ws.Cells.Item(113, 9).Value = 667518.5
$ws.Cells.Item(113, 10).Value = 2833
$ws.Cells.Item(113, 11).Value = 667518.5
$ws.Cells.Item(113, 12).Value = 2833
$ws.Cells.Item(113, 13).Value = -665348.5
$ws.Cells.Item(113, 14).Value = -7173
$ws.Cells.Item(119, 8).Value = 80000
$ws.Cells.Item(119, 10).Value = 80000
$ws.Cells.Item(119, 12).Value = 80000
$ws.Cells.Item(119, 14).Value = -89676
$ws.Cells.Item(134, 8).Value = 58996.2
$ws.Cells.Item(134, 10).Value = 58996.2
$ws.Cells.Item(134, 12).Value = 176988.6
$ws.Cells.Item(134, 14).Value = -182058.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2584.7778
$ws.Cells.Item(7, 9).Value = 2584.7778
$ws.Cells.Item(7, 11).Value = 2584.7778
$ws.Cells.Item(7, 13).Value = -2472.7778
$ws.Cells.Item(16, 8).Value = 1886.4546
$ws.Cells.Item(16, 9).Value = 1782
$ws.Cells.Item(16, 10).Value = 2069.25
$ws.Cells.Item(16, 11).Value = 1782
$ws.Cells.Item(16, 12).Value = 2069.25
$ws.Cells.Item(16, 13).Value = -1612
$ws.Cells.Item(16, 14).Value = -2409.25
$ws.Cells.Item(46, 8).Value = 10597.243
$ws.Cells.Item(46, 9).Value = 4179.5
$ws.Cells.Item(46, 11).Value = 4179.5
$ws.Cells.Item(46, 13).Value = -3991.5
$ws.Cells.Item(82, 8).Value = 493
$ws.Cells.Item(82, 9).Value = 249.66667
$ws.Cells.Item(82, 11).Value = 249.66667
$ws.Cells.Item(82, 13).Value = 111.33333
$ws.Cells.Item(85, 8).Value = 493
$ws.Cells.Item(85, 9).Value = 249.66667
$ws.Cells.Item(85, 11).Value = 249.66667
$ws.Cells.Item(85, 13).Value = 998.3333299999999
$ws.Cells.Item(93, 8).Value = 2831.0715
$ws.Cells.Item(93, 9).Value = 2923.125
$ws.Cells.Item(93, 11).Value = 2923.125
$ws.Cells.Item(93, 13).Value = -1675.125
$ws.Cells.Item(126, 8).Value = 2584.7778
$ws.Cells.Item(126, 9).Value = 2584.7778
$ws.Cells.Item(126, 11).Value = 7754.3334
$ws.Cells.Item(126, 13).Value = -5284.3334
$ws.Cells.Item(136, 8).Value = 5687.4736
$ws.Cells.Item(136, 9).Value = 4862.0713
$ws.Cells.Item(136, 10).Value = 7998.6
$ws.Cells.Item(136, 11).Value = 14586.2139
$ws.Cells.Item(136, 12).Value = 23995.8
$ws.Cells.Item(136, 13).Value = -12036.2139
$ws.Cells.Item(136, 14).Value = -29095.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = $null
$ws.Cells.Item(62, 8).Value = 4749.8184
$ws.Cells.Item(62, 9).Value = 4443.125
$ws.Cells.Item(62, 11).Value = 4443.125
$ws.Cells.Item(62, 13).Value = -3819.125
$ws.Cells.Item(65, 8).Value = 4749.8184
$ws.Cells.Item(65, 9).Value = 4443.125
$ws.Cells.Item(65, 11).Value = 22215.625
$ws.Cells.Item(65, 13).Value = -19095.625
$ws.Cells.Item(96, 8).Value = 1560.8572
$ws.Cells.Item(96, 9).Value = 1950
$ws.Cells.Item(96, 10).Value = 1405.2
$ws.Cells.Item(96, 11).Value = 1950
$ws.Cells.Item(96, 12).Value = 1405.2
$ws.Cells.Item(96, 13).Value = -577
$ws.Cells.Item(96, 14).Value = -4151.2
$ws.Cells.Item(100, 8).Value = 1007.3571
$ws.Cells.Item(100, 9).Value = 954
$ws.Cells.Item(100, 11).Value = 1908
$ws.Cells.Item(100, 13).Value = -1367
$ws.Cells.Item(113, 8).Value = 581.4167
$ws.Cells.Item(113, 9).Value = 577.7
$ws.Cells.Item(113, 11).Value = 1733.1
$ws.Cells.Item(113, 13).Value = 436.8999999999999
$ws.Cells.Item(126, 8).Value = 3221.5
$ws.Cells.Item(126, 9).Value = 3221.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 9664.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -7194.5
$ws.Cells.Item(126, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 4055.121
$ws.Cells.Item(132, 9).Value = 3710.6538
$ws.Cells.Item(132, 10).Value = 5334.5713
$ws.Cells.Item(132, 11).Value = 11131.9614
$ws.Cells.Item(132, 12).Value = 16003.7139
$ws.Cells.Item(132, 13).Value = -8601.9614
$ws.Cells.Item(132, 14).Value = -21063.7139
$ws.Cells.Item(136, 8).Value = 3700.2666
$ws.Cells.Item(136, 9).Value = 2448.6296
$ws.Cells.Item(136, 11).Value = 7345.888800000001
$ws.Cells.Item(136, 13).Value = -4795.888800000001
